# Refresh cryptos price/volume snapshot (scheduled GitHub Actions run)
# Columns: A=rank idx, B=Coin, C=Link, D=Price, E=Volume(1h)
# Rows 29-32 also got re-sorted (PancakeSwap/NEARProtocol/EthereumClassic/Aptos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold numeric-looking text (e.g. '68.246.70', '  +1.74%  ') that must stay
# plain text, not be reinterpreted as numbers/dates by Excel. A leading apostrophe
# forces text entry, exactly like typing it into Excel by hand.
function Set-TextCell($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}


# Row 2 - Bitcoin
Set-TextCell "D2" "68.222.53"
Set-TextCell "E2" "  +1.85%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.275.47"
Set-TextCell "E3" "  +0.50%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  -0.01%  "

# Row 5 - BNB
Set-TextCell "E5" "  +1.43%  "

# Row 6 - Solana
Set-TextCell "D6" "185.78"
Set-TextCell "E6" "  +4.99%  "

# Row 7 - USDC
Set-TextCell "E7" "  -0.06%  "

# Row 8 - XRP
Set-TextCell "D8" "0.602"
Set-TextCell "E8" "  -0.23%  "

# Row 9 - Dogecoin
Set-TextCell "E9" "  +4.13%  "

# Row 10 - Toncoin
Set-TextCell "D10" "6.74"
Set-TextCell "E10" "  -0.19%  "

# Row 11 - Cardano
Set-TextCell "D11" "0.418"
Set-TextCell "E11" "  +1.18%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextCell "D12" "3.843.58"
Set-TextCell "E12" "  +0.32%  "

# Row 13 - TRON
Set-TextCell "E13" "  +0.39%  "

# Row 14 - Avalanche
Set-TextCell "D14" "28.79"
Set-TextCell "E14" "  +2.37%  "

# Row 15 - WrappedBTC
Set-TextCell "D15" "68.249.82"
Set-TextCell "E15" "  +1.85%  "

# Row 16 - ShibaInu
Set-TextCell "E16" "  +3.08%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "3.272.64"
Set-TextCell "E17" "  +0.18%  "

# Row 18 - Polkadot
Set-TextCell "D18" "5.89"
Set-TextCell "E18" "  +1.02%  "

# Row 19 - Chainlink
Set-TextCell "D19" "13.66"
Set-TextCell "E19" "  +1.77%  "

# Row 20 - BitcoinCash
Set-TextCell "D20" "383.52"
Set-TextCell "E20" "  +3.06%  "

# Row 21 - Uniswap
Set-TextCell "D21" "7.74"
Set-TextCell "E21" "  +1.51%  "

# Row 22 - Dai
Set-TextCell "E22" "  +0.05%  "

# Row 23 - Litecoin
Set-TextCell "D23" "71.45"
Set-TextCell "E23" "  -0.07%  "

# Row 24 - Polygon
Set-TextCell "D24" "0.518"
Set-TextCell "E24" "  +0.89%  "

# Row 25 - PEPE
Set-TextCell "E25" "  +2.11%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextCell "D26" "9.93"
Set-TextCell "E26" "  +1.33%  "

# Row 27 - Kaspa
Set-TextCell "E27" "  +3.34%  "

# Row 28 - Binance-PegBSC-USD
Set-TextCell "D28" "1.00"
Set-TextCell "E28" "  +0.05%  "

# Row 29 - NEARProtocol (was PancakeSwap)
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D29" "5.81"
Set-TextCell "E29" "  +3.46%  "

# Row 30 - PancakeSwap (was NEARProtocol)
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D30" "2.00"
Set-TextCell "E30" "  +0.69%  "

# Row 31 - Aptos (was EthereumClassic)
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D31" "7.29"
Set-TextCell "E31" "  +6.72%  "

# Row 32 - EthereumClassic (was Aptos)
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D32" "22.99"
Set-TextCell "E32" "  +1.59%  "

# Row 33 - Fetch.AI
Set-TextCell "D33" "1.29"
Set-TextCell "E33" "  +1.37%  "

# Row 34 - USDe
Set-TextCell "D34" "0.998"
Set-TextCell "E34" "  +0.01%  "

# Row 35 - ImmutableX
Set-TextCell "E35" "  +3.07%  "

# Row 36 - Monero
Set-TextCell "D36" "162.70"
Set-TextCell "E36" "  -2.63%  "

# Row 37 - Stacks
Set-TextCell "D37" "1.88"
Set-TextCell "E37" "  +0.65%  "

# Row 38 - Mantle
Set-TextCell "D38" "0.839"
Set-TextCell "E38" "  -2.04%  "

# Row 39 - EnergySwap
Set-TextCell "D39" "26.85"
Set-TextCell "E39" "  -1.54%  "

# Row 40 - RenderToken
Set-TextCell "D40" "6.79"
Set-TextCell "E40" "  +4.70%  "

# Row 41 - Filecoin
Set-TextCell "E41" "  +6.11%  "

# Row 42 - dogwifhat
Set-TextCell "D42" "2.62"
Set-TextCell "E42" "  +1.58%  "

# Row 43 - Bittensor
Set-TextCell "D43" "350.34"
Set-TextCell "E43" "  +0.52%  "

# Row 44 - OKB
Set-TextCell "D44" "41.56"
Set-TextCell "E44" "  +2.54%  "

# Row 45 - InjectiveProtocol
Set-TextCell "D45" "25.66"
Set-TextCell "E45" "  +2.32%  "

# Row 46 - Hedera
Set-TextCell "D46" "0.0691"
Set-TextCell "E46" "  +2.12%  "

# Row 47 - Maker
Set-TextCell "D47" "2.657.87"
Set-TextCell "E47" "  -3.17%  "

# Row 48 - VeChain
Set-TextCell "E48" "  +1.78%  "

# Row 49 - Arweave
Set-TextCell "D49" "32.25"
Set-TextCell "E49" "  +5.80%  "

# Row 50 - ONDO
Set-TextCell "E50" "  +1.88%  "

# Row 51 - Stellar
Set-TextCell "E51" "  +0.15%  "
